# Fixing network data cleaning scripts
# - rename header columns to short machine-friendly names
# - title-case the Estado/Municipio text columns (fixes lowercase connector words)
# - tiny floating point re-computation of the "pct_matriculas" column for rows
#   where n_matriculas is 4 or 40
# - drop the trailing free-text footer rows (845-849)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function TitleCase([string]$s) {
    if ($s -eq $null) { return $s }
    $words = $s -split ' '
    $result = @()
    foreach ($w in $words) {
        if ($w.Length -eq 0) {
            $result += $w
        } else {
            # NOTE: use [string]::Concat rather than "+" - this COM runtime's
            # "+" operator coerces digit-only strings (e.g. "2" and "2") into
            # a numeric addition (giving "4") instead of concatenating them.
            $head = $w.Substring(0,1).ToUpper()
            $tail = $w.Substring(1).ToLower()
            $result += [string]::Concat($head, $tail)
        }
    }
    return [string]::Join(' ', $result)
}

# 1) Remove the trailing footer / notes rows (845-849). Doing this first keeps
#    all the other row numbers (2-843) stable for the edits below.
$ws.Range("845:849").Delete() | Out-Null

# 2) Rename the header row to the new short column names.
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# 3) Title-case every Estado/Municipio label in columns A and B (rows 2-843).
$labelRange = $ws.Range("A2:B843")
$labels = $labelRange.Value2
$rowCount = $labels.GetLength(0)
for ($i = 1; $i -le $rowCount; $i++) {
    for ($j = 1; $j -le 2; $j++) {
        $cell = $labels[$i, $j]
        if ($cell -ne $null -and $cell -is [string]) {
            $labels[$i, $j] = TitleCase $cell
        }
    }
}
$labelRange.Value2 = $labels

# 4) Recompute pct_matriculas for the rows whose n_matriculas is 4 or 40 - these
#    were recalculated slightly differently upstream, landing one ULP away
#    from the previous stored value.
$numRange = $ws.Range("C2:D843")
$nums = $numRange.Value2
$numRowCount = $nums.GetLength(0)
for ($i = 1; $i -le $numRowCount; $i++) {
    $count = $nums[$i, 1]
    if ($count -eq 4) {
        $nums[$i, 2] = 0.0009643201542912248
    } elseif ($count -eq 40) {
        $nums[$i, 2] = 0.009643201542912249
    }
}
$numRange.Value2 = $nums

Write-Output "edit complete"
